# refactor: change config after bha's discussion
#
# Re-layout the "DOMAIN" swim-lane on slide 1: grow the big background
# rectangle upward, nudge the DOMAIN label and the three existing
# accent-colored pill shapes up to make room, split the "Cleaning" pill's
# text into two runs, and add a brand-new pill ("Domain cleaning
# (discard values, fillna, etc.)") below "Engineering".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# 1) Big "DOMAIN" background rectangle -> grows upward and taller.
$domainBg = Get-ShapeByName $s "Rectangle 4"
$domainBg.Top = 159.13031496062993
$domainBg.Height = 316.1429921259843

# 2) "DOMAIN" label textbox -> shifts up slightly (and a 5-EMU nudge in X).
$domainLbl = Get-ShapeByName $s "TextBox 7"
$domainLbl.Left = 34.434251968503936
$domainLbl.Top = 173.41433070866142

# 3) "Machine learning" pill -> moves up.
$machineLearning = Get-ShapeByName $s "Rectangle 9"
$machineLearning.Top = 284.3404724409449

# 4) "Domain Intelligence" pill -> moves up.
$domainIntelligence = Get-ShapeByName $s "Rectangle 14"
$domainIntelligence.Top = 221.27818897637795

# 5) "Cleaning" pill -> re-worded into two runs.
$cleaning = Get-ShapeByName $s "Rectangle 15"
$cleaningText = $cleaning.TextFrame.TextRange
$cleaningText.Text = "Technical cleaning (encoding"
$cleaningText.InsertAfter(", misspelling errors, etc.)")

# 6) "Engineering" pill -> moves up.
$engineering = Get-ShapeByName $s "Rectangle 16"
$engineering.Top = 347.2328346456693

# 7) New pill below "Engineering": duplicate it to inherit the exact same
#    style/fill/line/size, then reposition and retext it.
$newPillRange = $engineering.Duplicate()
$newPill = $newPillRange.Item(1)
$newPill.Name = "Rectangle 12"
$newPill.Left = 43.825511811023624
$newPill.Top = 412.08614173228347

$newPillText = $newPill.TextFrame.TextRange
$newPillText.Text = "Domain cleaning (discard values, "
$newPillText.InsertAfter("fillna")
$newPillText.InsertAfter(", etc.)")
